# "update field beauty - MDLWL"
# The worksheet's "id" column (A2) holds a generated CA-xxxxxxxx token.
# This run regenerates/advances that token to a new value.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A2").Value = "CA-UACVGQ9R"
